# Fixes to rota text to address #116
#
# - Updates the Summer/Winter rota help text on the "setup" sheet
#   (previously "March to October" / "November to February",
#   now "April to September" / "October to March").
# - Leaves "setup" as the active sheet/selection (B17), and the
#   "welcome" sheet no longer tab-selected.

$wb = $excel.ActiveWorkbook

$wsSetup = $wb.Worksheets.Item("setup")

# Update the rota help strings with the corrected date ranges.
$wsSetup.Range("B17").Value = "☀️ Summer Rota runs from April to September"
$wsSetup.Range("B18").Value = "❄️ Winter Rota runs from October to March"

# Make "setup" the active sheet, with B17 selected.
$wsSetup.Activate()
$wsSetup.Range("B17").Select()
